# The workbook tracks daily Brócoli prices for "Feria Lagunitas de Puerto
# Montt". A new weekly observation was recorded and needs to be inserted
# as a new data row right after the existing header-adjacent rows, i.e. at
# worksheet row 117 (the data rows are otherwise sorted so this new
# observation slots in before the former row 117 and pushes every
# following row down by one, growing the used range from A1:R185 to
# A1:R186).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 117; everything from the old row 117
# through row 185 shifts down to rows 118-186 (values/styles move with
# their rows automatically).
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new observation.
$ws.Range("A117").Value = 4
$ws.Range("B117").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C117").Value = "Los Lagos"
$ws.Range("D117").Value = 44452
$ws.Range("E117").Value = 10
$ws.Range("F117").Value = 100112023
$ws.Range("G117").Value = "Brócoli"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Segunda"
$ws.Range("J117").Value = 500
$ws.Range("K117").Value = 1000
$ws.Range("L117").Value = 1000
$ws.Range("M117").Value = 1000
$ws.Range("N117").Value = "`$/unidad"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 1000
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"

# Make sure the new row's date cell carries the same date number format
# ("s=2") used by every other row's Fecha column.
$ws.Range("D117").NumberFormat = $ws.Range("D118").NumberFormat
